$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 'Numbers to 20'),
    @(2, 'Days of the week, months of the year, parts of the day'),
    @(3, 'Colors'),
    @(4, 'Simple conversational phrases and essential words'),
    @(5, 'Question words, articles, other basic words'),
    @(6, 'Conjugations of the most basic verbs'),
    @(7, 'Verb infinitives'),
    @(8, 'Conjugations of regular verbs'),
    @(10, 'Food'),
    @(12, 'Cooking and serving food'),
    @(14, 'Conjugations of some irregular verbs'),
    @(16, 'Time of day'),
    @(18, 'Numbers 20 to 100, 200, 1000, 2000'),
    @(20, 'People'),
    @(22, 'Relatives'),
    @(24, 'Things from daily life, things you have with you'),
    @(26, 'Miscellaneous 1'),
    @(28, 'Possesives'),
    @(30, 'Conversational phrases'),
    @(32, 'Verb infinitives'),
    @(34, 'Parts of a house'),
    @(36, 'Things around the house'),
    @(38, 'Miscellaneous 2'),
    @(40, 'Verb examples'),
    @(42, 'Animals'),
    @(44, 'Weather'),
    @(46, 'Miscellaneous 3'),
    @(48, 'School, classroom'),
    @(50, 'Verb infinitives'),
    @(52, 'Miscellaneous 4'),
    @(54, 'Words from books'),
    @(56, 'Clothing'),
    @(58, 'Things in your neighborhood'),
    @(60, 'Places around town, getting around town'),
    @(61, 'Travel'),
    @(62, 'Shopping'),
    @(63, 'Miscellaneous 5'),
    @(64, 'Restaurants'),
    @(66, 'Days of the year, holidays'),
    @(67, 'Entertainment'),
    @(68, 'Parts of the body'),
    @(69, 'Doctors, medicine, health'),
    @(70, 'Finding a place to live, apartments, renting'),
    @(71, 'Personality traits'),
    @(72, 'Conversational phrases'),
    @(73, 'Things in the country'),
    @(74, 'Cars, driving'),
    @(75, 'Post office, mail'),
    @(76, 'Plurals'),
    @(77, 'Ordinal numbers'),
    @(78, 'Food plurals'),
    @(79, 'Clothing plurals'),
    @(81, 'Accusative'),
    @(82, 'Past tense'),
    @(83, 'Present and past participles of verbs'),
    @(84, 'Verb infinitives'),
    @(86, 'Present past'),
    @(87, 'University classes'),
    @(89, 'Sports'),
    @(90, 'Government'),
    @(91, 'Relatives (advanced)'),
    @(92, 'Food (advanced)'),
    @(93, 'Places and things in town (advanced)'),
    @(94, 'Things from around the house (advanced)'),
    @(95, 'Miscellaneous and maybe not very common'),
    @(96, 'Miscellaneous and maybe not very common'),
    @(97, 'Contracts'),
    @(98, 'Sentences'),
    @(99, 'Professions, business'),
)

# Seed brand-new shared strings in the same order they appear in the
# target file's sharedStrings.xml (new strings are appended in order of
# first write, so "Words from books" (row 31) must be written before
# "Verb examples" (row 24)).
$ws.Cells.Item(31, 2).Value = "Words from books"
$ws.Cells.Item(24, 2).Value = "Verb examples"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A68"))
$ws.Sort.SetRange($ws.Range("A1:B68"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

$ws.Range("A25").Select()

Write-Host "done"